$d = $word.ActiveDocument
$r = $d.Content
$r.Collapse(0)
$r.InsertAfter(" with e")
$r.Collapse(0)
$r.InsertAfter("xtension")
